# Add the new "altar" entity to the collision table on Sheet1.
# This mirrors the existing "player" row/column: altar does not collide
# with player ("▬"), but collides with everything else ("x"), including
# itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column header (K1) and new row header (A11) both use the new shared
# string "altar".
$ws.Range("K1").Value = "altar"
$ws.Range("A11").Value = "altar"

# New column K, rows 2-10 (collision of each existing entity vs altar).
# Row 2 is "player" -> altar does not collide with player.
$ws.Range("K2").Value = "▬"
$ws.Range("K3").Value = "x"
$ws.Range("K4").Value = "x"
$ws.Range("K5").Value = "x"
$ws.Range("K6").Value = "x"
$ws.Range("K7").Value = "x"
$ws.Range("K8").Value = "x"
$ws.Range("K9").Value = "x"
$ws.Range("K10").Value = "x"

# New row 11 (altar vs each existing entity, including the new altar column).
# Column B is "player" -> altar does not collide with player.
$ws.Range("B11").Value = "▬"
$ws.Range("C11").Value = "x"
$ws.Range("D11").Value = "x"
$ws.Range("E11").Value = "x"
$ws.Range("F11").Value = "x"
$ws.Range("G11").Value = "x"
$ws.Range("H11").Value = "x"
$ws.Range("I11").Value = "x"
$ws.Range("J11").Value = "x"
$ws.Range("K11").Value = "x"

# Update the active cell / selection to match the post-edit cursor position.
$ws.Range("I16").Select()
